# Uppdatera uppgift och tester funkar
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tidsredovisning")

$tbl = $ws.ListObjects.Item("Tabell2")

# Remove the totals row first so new data rows can be appended to the
# bottom of the table without the SUBTOTAL row getting in the way.
$tbl.TotalsRowShown = $false

# New rows of time-tracking data to append to the "Tabell2" table.
$newRows = @(
    @(45315, 2, "PHP Spara ny post"),
    @(45317, 1, "PHP Spara ny post tester"),
    @(45317, 1, "PHP Kontrollera indata tester"),
    @(45317, 2, "PHP Hämta enskild uppgift + test")
)

foreach ($row in $newRows) {
    $listRow = $tbl.ListRows.Add()
    $r = $listRow.Range
    $r.Cells.Item(1, 1).Value = $row[0]
    $r.Cells.Item(1, 1).NumberFormat = "m/d/yyyy"
    $r.Cells.Item(1, 2).Value = $row[1]
    $r.Cells.Item(1, 3).Value = $row[2]
}

# Restore the totals row with the SUBTOTAL formula.
$tbl.TotalsRowShown = $true

# Update the view to match the edited selection/scroll position.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("F15").Select()
